$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a date-looking literal string ("YYYY-MM-DD") into a cell
# without Excel's autodetection turning it into a date serial number.
# Trick: force the cell to Text format ("@") before assigning the value,
# then reset the cell style back to "Normal" so no lingering style index
# is left on the cell (matches the source file, where these cells carry
# no explicit style).
# ---------------------------------------------------------------------------
function Set-LiteralText {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ===========================================================================
# Sheet "Summary"
# ===========================================================================
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1401.37
$wsSummary.Range("B4").Value = 1.16
$wsSummary.Range("B6").Value = 129
$wsSummary.Range("B8").Value = 50
$wsSummary.Range("B9").Value = 42.64

# ===========================================================================
# Sheet "Strategy Status" - MarketMaking row (row 5)
# ===========================================================================
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 101.37
$wsStatus.Range("D5").Value = 96
$wsStatus.Range("E5").Value = 1.05
$wsStatus.Range("F5").Value = 1.37
$wsStatus.Range("G5").Value = 42.71

# ===========================================================================
# Sheet "All Trades"
# ===========================================================================
$wsAll = $wb.Worksheets.Item("All Trades")

# Trade #129 (row 130) gets closed out early.
$wsAll.Range("G130").Value = 0.074143
$wsAll.Range("H130").Value = "CLOSED"
$wsAll.Range("I130").Value = -17.6192
$wsAll.Range("J130").Value = -0.02
$wsAll.Range("K130").Value = 101.37
$wsAll.Range("L130").Value = "early_exit"
$wsAll.Range("M130").Value = 0.15

# New trade #162 opened - appended as row 163.
$wsAll.Range("A163").Value = 162
Set-LiteralText $wsAll.Range("B163") "2026-02-17"
$wsAll.Range("C163").Value = "21:27:36"
$wsAll.Range("D163").Value = "MarketMaking"
$wsAll.Range("E163").Value = "UP"
$wsAll.Range("F163").Value = 0.09
$wsAll.Range("H163").Value = "OPEN"
$wsAll.Range("I163").Value = 0
$wsAll.Range("J163").Value = 0
$wsAll.Range("K163").Value = 101.3900330787957
$wsAll.Range("M163").Value = 0
$wsAll.Range("N163").Value = 0
$wsAll.Range("O163").Value = 0
$wsAll.Range("P163").Value = 0.6
$wsAll.Range("Q163").Value = "Normal spread capture: 19600 bps"

# ===========================================================================
# Sheet "MarketMaking" (per-strategy trade log)
# ===========================================================================
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Trade #129 (row 97 on this sheet) gets closed out early.
$wsMM.Range("G97").Value = 0.074143
$wsMM.Range("H97").Value = "CLOSED"
$wsMM.Range("I97").Value = -17.6192
$wsMM.Range("J97").Value = -0.02
$wsMM.Range("K97").Value = 101.37
$wsMM.Range("P97").Value = "early_exit"
$wsMM.Range("Q97").Value = 0.15

# New trade #162 opened - appended as row 130 on this sheet.
$wsMM.Range("A130").Value = 162
Set-LiteralText $wsMM.Range("B130") "2026-02-17"
$wsMM.Range("C130").Value = "21:27:36"
$wsMM.Range("D130").Value = "MarketMaking"
$wsMM.Range("E130").Value = "UP"
$wsMM.Range("F130").Value = 0.09
$wsMM.Range("H130").Value = "OPEN"
$wsMM.Range("I130").Value = 0
$wsMM.Range("J130").Value = 0
$wsMM.Range("K130").Value = 101.3900330787957
$wsMM.Range("L130").Value = 0
$wsMM.Range("M130").Value = 0
$wsMM.Range("N130").Value = 0.6
$wsMM.Range("O130").Value = "Normal spread capture: 19600 bps"
$wsMM.Range("Q130").Value = 0
